$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that can look like numbers (e.g. "1.019").
# Force them to be stored as text (matching the inline string source data) by
# temporarily applying a text number format, then restore the default style so
# the cell keeps its original (unstyled) appearance.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.877.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.019"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4323"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07434"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8867"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.864.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.782"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07119"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.76%  "
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009056"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("E19").Value = "  -1.70%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.926.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.287"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.096.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.027"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.433"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08994"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.243"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7782"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.602"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.937"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.018"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05329"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01972"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.891"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5210"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.013"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1684"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.785"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.723"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4761"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.017"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("E51").Value = "  +0.75%  "
